$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "FORG-..." external reference code used in column C (rows 2-5)
$ws.Range("C2:C5").Value = "FORG-7330-1645-9608"

# Update the AWS-account-like GUID used in column K (rows 2-5)
$ws.Range("K2:K5").Value = "06cb0fb3-9a1d-4644-815e-f6d13b8faa18"

# Update the GCP-account GUID used in column M rows 2-3
$ws.Range("M2:M3").Value = "a65ed389-b74d-4fb1-b2e7-06298fadc1a6"

# Update the GCP-account GUID used in column M rows 4-5
$ws.Range("M4:M5").Value = "386ecbe0-6b65-46a4-8b81-48f6b38a088c"
